$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "304.20"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.81%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.69"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.17%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.173"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.22%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07480"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.59%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.405"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "39.34%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.014"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.96%"

$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9147"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.10%"

$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1731"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.27%"

$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07688"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "4.48%"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08180"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.95%"

$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03027"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.02%"

$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09937"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.38%"

$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001511"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.57%"

$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.006150"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.07%"

$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.496"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.33%"

$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.867"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.87%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.237"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.72%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3262"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.89%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1329"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.38%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.657"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.37%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04617"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.76%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.1565"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.02%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001262"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3.82%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004529"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.36%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001299"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-7.14%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002740"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "48.82%"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "5.42%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04533"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.33%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007384"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "5.37%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1362"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "2.64%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002168"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "5.35%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01089"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-14.68%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006391"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "6.08%"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "15.31%"
